$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 new rows before the old blank separator row (58), pushing
# the separator row and the total row down to 62 and 63.
$ws.Rows("58:61").Insert(-4121)

# --- Row 58: Skype Session ---
$ws.Range("A58").Value = "Skype Session"
$ws.Range("B58").Value = 43743
$ws.Range("C58").Value = 0.5
$ws.Range("D58").Value = 25
$ws.Range("E58").Formula = "=D58*C58"

# --- Row 59: Attempts at Image Orientation ---
$ws.Range("A59").Value = "Attempts at Image Orientation"
$ws.Range("B59").Value = 43744
$ws.Range("C59").Value = 1
$ws.Range("D59").Value = 25
$ws.Range("E59").Formula = "=D59*C59"

# --- Row 60: Skype Session ---
$ws.Range("A60").Value = "Skype Session"
$ws.Range("B60").Value = 43753
$ws.Range("C60").Value = 0.5
$ws.Range("D60").Value = 25
$ws.Range("E60").Formula = "=D60*C60"

# --- Row 61: Changing Price display ---
$ws.Range("A61").Value = "Changing Price display"
$ws.Range("B61").Value = 43753
$ws.Range("C61").Value = 0.5
$ws.Range("D61").Value = 25
$ws.Range("E61").Formula = "=D61*C61"

# --- Row 62 (blank separator row, pushed down from row 58) ---
$ws.Range("B61").Copy()
$ws.Range("B62").PasteSpecial(-4122)
$ws.Range("E61").Copy()
$ws.Range("E62").PasteSpecial(-4122)

# --- Row 63 (total row, pushed down from row 59) ---
$ws.Range("C63").Formula = "=SUM(C53:C61)"
$ws.Range("E63").Formula = "=SUM(E53:E61)"

# Fix up the active selection to match the final cursor position
$null = $ws.Range("D62").Select()

Write-Host "done"
